# "added live demo test case" - refresh the canned Funding / Expense demo
# dates two years into the future, and move the active selection.

$wb = $excel.ActiveWorkbook
$wsFunding = $wb.Worksheets.Item("Funding")
$wsExpense = $wb.Worksheets.Item("Expense")

# --- Funding sheet: push every Start/End date out by 2 years ---------------
$wsFunding.Range("C2").Value = 46388
$wsFunding.Range("D2").Value = 46568
$wsFunding.Range("C3").Value = 46569
$wsFunding.Range("D3").Value = 46752
$wsFunding.Range("C4").Value = 46388
$wsFunding.Range("D4").Value = 46752
$wsFunding.Range("C5").Value = 46388
$wsFunding.Range("D5").Value = 46752
$wsFunding.Range("C6").Value = 46388
$wsFunding.Range("D6").Value = 46752

# --- Expense sheet: push every date out by 2 years --------------------------
$wsExpense.Range("E2").Value = 46461
$wsExpense.Range("E3").Value = 46433
$wsExpense.Range("E4").Value = 46433
$wsExpense.Range("E5").Value = 46492
$wsExpense.Range("E6").Value = 46619
$wsExpense.Range("E7").Value = 46553
$wsExpense.Range("E8").Value = 46753
$wsExpense.Range("E9").Value = 46037

# --- Update selection on each sheet, then make Funding the active tab ------
$wsExpense.Range("D13").Select()
$wsFunding.Range("C15").Select()
$wsFunding.Activate()
